# CRIAR_PEDIDOS.vbs refactor: update the purchase-order / SAP order numbers
# that were filled in once the orders P2/P3 (and their related Q2/Q3
# contract items) came back from SAP, then leave the sheet positioned/
# selected the way the workbook was left after the save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: order P2207 -> P2209 (corrected order number) and Q2 now has the
# matching contract/item number that came back from SAP.
$ws.Range("P2").Value = 4600244209
$ws.Range("Q2").Value = 4503342003

# Row 3: P3/Q3 were still blank - they now got their SAP order numbers too.
$ws.Range("P3").Value = 4600244210
$ws.Range("Q3").Value = 4503342004

# Reposition the view/selection the way it was left when the workbook was
# saved: scrolled so column H is the left-most visible column, with P3
# selected/active.
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 1
$null = $ws.Range("P3").Select()

$null = $wb.Save()
